# Improved main infographic in chapter 3
#
# Slide 1 contains the integrated-experimental-strategy infographic.
# Two vertical "spine" labels on the left-hand side of the figure
# ("RESEARCH QUESTIONS" and "RESEARCH OBJECTIVES") are repositioned /
# resized so the figure reads better:
#   - TextBox 88 (RESEARCH QUESTIONS) moves down and gets a touch shorter.
#   - TextBox 97 (RESEARCH OBJECTIVES) moves up slightly and gets taller.
# Left/Width stay the same for both; only Top/Height change.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# PowerPoint's Shape.Top/.Height are expressed in points (1 pt = 12700 EMU)
# and are stored internally as single-precision floats, so the literals
# below are chosen to land exactly on the target EMU values after that
# float32 round-trip.

$questions = $s.Shapes.Item(32)   # "TextBox 88" / RESEARCH QUESTIONS
$questions.Top    = 75.6611099243164    # -> 960896 EMU (was 831652 EMU)
$questions.Height = 89.86071014404297   # -> 1141231 EMU (was 1239149 EMU)

$objectives = $s.Shapes.Item(41)  # "TextBox 97" / RESEARCH OBJECTIVES
$objectives.Top    = 167.11976623535156 # -> 2122421 EMU (was 2138663 EMU)
$objectives.Height = 98.84968566894531  # -> 1255391 EMU (was 1239149 EMU)
